$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving a purely-numeric-looking Price value need the cell
# pre-formatted as Text so Excel keeps them as strings (preserving
# trailing zeros / exact formatting) instead of coercing to a Number.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated values (Coin/Link/Price/Volume) cell by cell.
$ws.Range("D2").Value = "28.472.09"
$ws.Range("E2").Value = "  +2.42%  "
$ws.Range("D3").Value = "1.865.19"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("D4").Value = "1.014"
$ws.Range("E4").Value = "  +0.71%  "
$ws.Range("D5").Value = "335.79"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("E6").Value = "  +1.11%  "
$ws.Range("D7").Value = "0.4582"
$ws.Range("E7").Value = "  -1.87%  "
$ws.Range("D8").Value = "0.3953"
$ws.Range("E8").Value = "  +2.39%  "
$ws.Range("D9").Value = "47.72"
$ws.Range("E9").Value = "  +1.78%  "
$ws.Range("D10").Value = "0.07890"
$ws.Range("E10").Value = "  -0.41%  "
$ws.Range("D11").Value = "0.9850"
$ws.Range("E11").Value = "  +1.65%  "
$ws.Range("D12").Value = "21.53"
$ws.Range("E12").Value = "  +1.01%  "
$ws.Range("D13").Value = "1.894.84"
$ws.Range("E13").Value = "  +2.38%  "
$ws.Range("D14").Value = "5.869"
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("D15").Value = "7.024"
$ws.Range("E15").Value = "  -1.50%  "
$ws.Range("D16").Value = "1.013"
$ws.Range("E16").Value = "  +0.57%  "
$ws.Range("D17").Value = "88.26"
$ws.Range("E17").Value = "  -2.90%  "
$ws.Range("D18").Value = "0.06595"
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("D19").Value = "0.00001030"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").Value = "17.14"
$ws.Range("E20").Value = "  -0.86%  "
$ws.Range("D21").Value = "1.010"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").Value = "28.485.31"
$ws.Range("E22").Value = "  +2.49%  "
$ws.Range("D23").Value = "5.381"
$ws.Range("E23").Value = "  +0.65%  "
$ws.Range("D24").Value = "10.75"
$ws.Range("E24").Value = "  -0.45%  "
$ws.Range("E25").Value = "  -1.03%  "
$ws.Range("D26").Value = "2.128.37"
$ws.Range("E26").Value = "  +2.67%  "
$ws.Range("D27").Value = "157.45"
$ws.Range("E27").Value = "  -1.11%  "
$ws.Range("D28").Value = "19.37"
$ws.Range("E28").Value = "  -0.45%  "
$ws.Range("D29").Value = "2.062"
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("D30").Value = "5.338"
$ws.Range("E30").Value = "  -0.95%  "
$ws.Range("D31").Value = "117.34"
$ws.Range("E31").Value = "  -1.04%  "
$ws.Range("D32").Value = "0.9483"
$ws.Range("E32").Value = "  +0.76%  "
$ws.Range("D33").Value = "0.09368"
$ws.Range("D34").Value = "3.605"
$ws.Range("D35").Value = "1.394"
$ws.Range("E35").Value = "  +5.03%  "
$ws.Range("D36").Value = "5.254"
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("D37").Value = "0.06035"
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("D38").Value = "0.02213"
$ws.Range("E38").Value = "  -0.18%  "
$ws.Range("D39").Value = "8.206"
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("D40").Value = "1.155"
$ws.Range("E40").Value = "  -0.54%  "
$ws.Range("D41").Value = "1.014"
$ws.Range("E41").Value = "  +0.81%  "
$ws.Range("D42").Value = "0.5787"
$ws.Range("E42").Value = "  -0.45%  "
$ws.Range("E43").Value = "  +0.26%  "
$ws.Range("D44").Value = "0.1813"
$ws.Range("E44").Value = "  -1.91%  "
$ws.Range("D45").Value = "1.232"
$ws.Range("E45").Value = "  -4.04%  "
$ws.Range("D46").Value = "2.299"
$ws.Range("E46").Value = "  +28.40%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "0.5459"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "11.84"
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("D49").Value = "0.07173"
$ws.Range("E49").Value = "  +4.93%  "
$ws.Range("D50").Value = "1.887"
$ws.Range("E50").Value = "  -2.59%  "
$ws.Range("D51").Value = "110.37"
$ws.Range("E51").Value = "  -0.60%  "
